# fix: router and db
# Adds 5 new transportadora/delivery rows (7-11) to Sheet1, updates the
# "Estado" values on rows 2-6 from the stale RJ/SP placeholders to their
# proper text, and turns the leftover shared-string placeholder
# cells in D5/D6/G5/G6 back into plain numeric 0 / phone-number cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rows 2-4: Estado column just re-affirmed (RJ/SP/RJ) ---
$ws.Range("F2").Value = "RJ"
$ws.Range("F3").Value = "SP"
$ws.Range("F4").Value = "RJ"

# --- rows 5-6: Preco de km/telefone placeholders become real numbers ---
$ws.Range("D5").Value = 0
$ws.Range("F5").Value = "RJ"
$ws.Range("G5").Value = 21888888888

$ws.Range("D6").Value = 0
$ws.Range("F6").Value = "RJ"
$ws.Range("G6").Value = 21888888888

# --- row 7: granja transportes ---
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "granja transportes"
$ws.Range("C7").Value = "'15"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = "granja@transportes"
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = "CE"
$ws.Range("G7").Value = "'67999994444"
$ws.Range("G7").Style = "Normal"

# --- row 8: transportadora caina ---
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "transportadora caina"
$ws.Range("C8").Value = "'231"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = "caina@email"
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = "ES"
$ws.Range("G8").Value = "'11111111111"
$ws.Range("G8").Style = "Normal"

# --- row 9: transportadora stefany ---
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "transportadora stefany"
$ws.Range("C9").Value = "'123"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = "stefany@email"
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = "AM"
$ws.Range("G9").Value = "'123123"
$ws.Range("G9").Style = "Normal"

# --- row 10: transporte delivery ---
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "transporte delivery"
$ws.Range("C10").Value = "'111"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "transporte@delivery"
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = "MT"
$ws.Range("G10").Value = "'123123123123"
$ws.Range("G10").Style = "Normal"

# --- row 11: delivery ---
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "delivery"
$ws.Range("C11").Value = "'1"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "delivery@email"
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = "PE"
$ws.Range("G11").Value = "'123123321321"
$ws.Range("G11").Style = "Normal"
